# SARAALERT-1260: Allow vaccine table to be populated on import
# Adds "Vaccine 1 ..." and "Vaccine 2 ..." columns (CY:DH) to the Monitorees
# export/import template, with sample data rows for the first few monitorees.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): new Vaccine 1 / Vaccine 2 columns
# ---------------------------------------------------------------------------
$ws.Range("CY1").Value2 = "Vaccine 1 Group Name"
$ws.Range("CZ1").Value2 = "Vaccine 1 Product Name"
$ws.Range("DA1").NumberFormat = "@"
$ws.Range("DA1").Value2 = "Vaccine 1 Administration Date"
$ws.Range("DB1").Value2 = "Vaccine 1 Dose Number"
$ws.Range("DC1").Value2 = "Vaccine 1 Notes"
$ws.Range("DD1").Value2 = "Vaccine 2 Group Name"
$ws.Range("DE1").Value2 = "Vaccine 2 Product Name"
$ws.Range("DF1").NumberFormat = "@"
$ws.Range("DF1").Value2 = "Vaccine 2 Administration Date"
$ws.Range("DG1").Value2 = "Vaccine 2 Dose Number"
$ws.Range("DH1").Value2 = "Vaccine 2 Notes"

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$ws.Range("CY2").Value2 = "COVID-19"
$ws.Range("CZ2").Value2 = "Moderna COVID-19 Vaccine"
$ws.Range("DA2").NumberFormat = "@"
$ws.Range("DA2").Value2 = "2020-06-01"
$ws.Range("DB2").Value2 = 1
$ws.Range("DC2").Value2 = "notes 1"
$ws.Range("DD2").Value2 = "COVID-19"
$ws.Range("DE2").Value2 = "Moderna COVID-19 Vaccine"
$ws.Range("DF2").NumberFormat = "@"
$ws.Range("DF2").Value2 = "2020-06-20"
$ws.Range("DG2").Value2 = 2
$ws.Range("DH2").Value2 = "notes 2"

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Range("CY3").Value2 = "COVID-19"
$ws.Range("CZ3").Value2 = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DA3").NumberFormat = "@"
$ws.Range("DA3").Value2 = "2020-06-02"
$ws.Range("DB3").Value2 = 1
$ws.Range("DD3").Value2 = "COVID-19"
$ws.Range("DE3").Value2 = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("DF3").NumberFormat = "@"
$ws.Range("DF3").Value2 = "2020-06-21"
$ws.Range("DG3").Value2 = 2

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("CY4").Value2 = "COVID-19"
$ws.Range("CZ4").Value2 = "Unknown"
$ws.Range("DA4").NumberFormat = "@"
$ws.Range("DA4").Value2 = "2020-06-04"
$ws.Range("DB4").Value2 = 1
$ws.Range("DD4").Value2 = "COVID-19"
$ws.Range("DE4").Value2 = "Unknown"
$ws.Range("DF4").NumberFormat = "@"
$ws.Range("DF4").Value2 = "2020-06-22"
$ws.Range("DG4").Value2 = 2

# ---------------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------------
$ws.Range("CY5").Value2 = "COVID-19"
$ws.Range("CZ5").Value2 = "Moderna COVID-19 Vaccine"
$ws.Range("DA5").NumberFormat = "@"
$ws.Range("DA5").Value2 = "2020-06-01"
$ws.Range("DB5").Value2 = 1

# ---------------------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------------------
$ws.Range("CY6").Value2 = "COVID-19"
$ws.Range("CZ6").Value2 = "Janssen (J&J) COVID-19 Vaccine"
$ws.Range("DA6").NumberFormat = "@"
$ws.Range("DA6").Value2 = "2020-06-03"
$ws.Range("DB6").Value2 = 1

# ---------------------------------------------------------------------------
# Row 7
# ---------------------------------------------------------------------------
$ws.Range("CY7").Value2 = "COVID-19"
$ws.Range("CZ7").Value2 = "Unknown"
$ws.Range("DA7").NumberFormat = "@"
$ws.Range("DA7").Value2 = "2020-06-02"
$ws.Range("DB7").Value2 = 1

# ---------------------------------------------------------------------------
# Column widths for the new columns (best-fit, matching the widths Excel
# computed for these headers/values in the authored workbook).
# ---------------------------------------------------------------------------
$ws.Columns.Item(103).ColumnWidth = 19.571428571428573   # CY  -> stored 20.33203125
$ws.Columns.Item(104).ColumnWidth = 30.285714285714285   # CZ  -> stored 31
$ws.Columns.Item(105).ColumnWidth = 25                   # DA  -> stored 25.6640625
$ws.Columns.Item(106).ColumnWidth = 20.428571428571427   # DB  -> stored 21.1640625
$ws.Columns.Item(107).ColumnWidth = 13.857142857142858   # DC  -> stored 14.5
$ws.Columns.Item(108).ColumnWidth = 19.571428571428573   # DD  -> stored 20.33203125
$ws.Columns.Item(109).ColumnWidth = 30.285714285714285   # DE  -> stored 31
$ws.Columns.Item(110).ColumnWidth = 25                   # DF  -> stored 25.6640625
$ws.Columns.Item(111).ColumnWidth = 20.428571428571427   # DG  -> stored 21.1640625
$ws.Columns.Item(112).ColumnWidth = 13.857142857142858   # DH  -> stored 14.5

# ---------------------------------------------------------------------------
# Reset the view: scroll back to the top-left and select A1 (clears the old
# topLeftCell="CO1" / selection="CU9" saved in the sheet view).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
